$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -8.016
$ws.Range("C9").Value = -11.002
$ws.Range("D9").Value = -7.356
$ws.Range("D11").Value = -7.323
$ws.Range("C18").Value = -11.888
$ws.Range("C20").Value = -12.278
$ws.Range("D23").Value = -8.266
$ws.Range("D24").Value = -7.100999999999999
$ws.Range("D26").Value = -7.725
$ws.Range("C27").Value = -12.135
$ws.Range("D34").Value = -7.526999999999999
$ws.Range("C35").Value = -12.092
$ws.Range("D35").Value = -8.089000000000002
$ws.Range("D48").Value = -7.781000000000001
$ws.Range("D49").Value = -7.867
$ws.Range("D52").Value = -7.926
$ws.Range("D66").Value = -7.577
$ws.Range("D67").Value = -7.733
$ws.Range("C69").Value = -11.105
$ws.Range("C76").Value = -12.969
$ws.Range("C78").Value = -12.296
$ws.Range("D78").Value = -7.218999999999999
$ws.Range("D80").Value = -8.316999999999998
$ws.Range("C82").Value = -11.549
$ws.Range("C83").Value = -13.105
$ws.Range("C93").Value = -10.818
$ws.Range("D99").Value = -7.968000000000001
$ws.Range("D104").Value = -7.647
